$wb = $excel.ActiveWorkbook

# "Ohjelmointi" sheet holds the command-parameter table; B7/B8 share the
# "target temperature" description string. Extend it to mention the new
# one-decimal precision now available from the analog thermistor read.
$ws = $wb.Worksheets.Item("Ohjelmointi")
$ws.Range("B7").Value = "tavoitelämpötila(0 - 25,4)"
$ws.Range("B8").Value = "tavoitelämpötila(0 - 25,4)"

# Move the active selection to D6 (matches the author's last click before save)
$ws.Activate()
$ws.Range("D6").Select()
